# Generate Report for Archive
#
# The localization status report was regenerated:
#  1) Rows that were previously "Ready for handoff" are now "In Translation"
#     (this text lives in the Overview sheet's per-locale status columns
#     E/F, and in the Status column (C) of each per-locale detail sheet).
#  2) The Status column(s) that used to hold the longer "Ready for handoff"
#     text are narrowed now that the shorter "In Translation" text fits.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Replace every occurrence of the old status text with the new one ---
foreach ($ws in @($overview, $zhcn, $dede)) {
    $used = $ws.UsedRange
    foreach ($row in 1..$used.Rows.Count) {
        foreach ($col in 1..$used.Columns.Count) {
            $cell = $ws.Cells.Item($row, $col)
            $cellText = [string]$cell.Text
            if ($cellText -eq $oldStatus) {
                $cell.Value = $newStatus
            }
        }
    }
}

# --- Narrow the Status column(s) now that their content is shorter ---
$overview.Range("E1:F1").ColumnWidth = 12.5
$zhcn.Range("C1").ColumnWidth = 12.5
$dede.Range("C1").ColumnWidth = 12.5
